$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Apply the new "explicit Calibri 10" font to the existing tipo_maquina
#    (column B) cells, rows 1-5 (header + the four original data rows).
#    Re-asserting the font name (even though it is already Calibri) forces
#    Excel to register an explicit (non theme-linked) font record, matching
#    the new font that shows up in the diff.
# ---------------------------------------------------------------------------
$ws.Range("B1:B5").Font.Name = "Calibri"

# ---------------------------------------------------------------------------
# 2. Fix up the existing data in row 2 (Escavadeira Hidraulica):
#      E2: Preventiva        -> Sem ocorrencias
#      F2: Sem ocorrências   -> Sem ocorrencias
#      I2: "3500.0" (text)   -> 0 (number)
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = "Sem ocorrencias"
$ws.Range("F2").Value = "Sem ocorrencias"
$ws.Range("I2").Value = 0

# ---------------------------------------------------------------------------
# 3. Append three new rows (5, 6, 7 -> sheet rows 6, 7, 8) of machine data.
#    Using Insert() (rather than simply writing into blank rows) makes the
#    new rows inherit the number formatting of the row above them, so the
#    H column picks up the same date style used by H2:H5 automatically.
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).Insert() | Out-Null
$ws.Rows.Item(7).Insert() | Out-Null
$ws.Rows.Item(8).Insert() | Out-Null

$ws.Rows.Item(6).RowHeight = 15.75
$ws.Rows.Item(7).RowHeight = 15.75
$ws.Rows.Item(8).RowHeight = 15.75

# Row 6 - caminhões de mineração
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "caminhões de mineração"
$ws.Range("C6").Value = "Inativa"
$ws.Range("D6").Value = "Ana Souza"
$ws.Range("E6").Value = "Corretiva"
$ws.Range("F6").Value = "Muito Alta"
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 45589

# Row 7 - perfuratrizes (first occurrence)
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "perfuratrizes "
$ws.Range("C7").Value = "Em manutenção"
$ws.Range("D7").Value = "Ana Souza"
$ws.Range("E7").Value = "Corretiva"
$ws.Range("F7").Value = "Alta"
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 45658

# Row 8 - perfuratrizes (second occurrence)
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "perfuratrizes "
$ws.Range("C8").Value = "Ativa"
$ws.Range("D8").Value = "Ana Souza"
$ws.Range("E8").Value = "Sem ocorrencias"
$ws.Range("F8").Value = "Sem ocorrencias"
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 45940
$ws.Range("I8").Value = 0

# Give the "tipo_maquina" cells of the three new rows their own, larger font
# (matches the bigger font used for the new rows in the diff).
$ws.Range("B6:B8").Font.Size = 12

# ---------------------------------------------------------------------------
# 4. I6 and I7 need to hold the numbers "7000.0" / "6000.0" as literal TEXT
#    (matching custo_total_manutencao for the other rows), but assigning a
#    numeric-looking string straight to .Value causes Excel to coerce it to
#    a real number. Work around this by formatting a scratch cell as Text,
#    writing the value there, copying it, and pasting values-only into the
#    destination cell (which keeps the General/default style on the target
#    cell while still storing the content as text).
# ---------------------------------------------------------------------------
$ws.Range("K1").NumberFormat = "@"

$ws.Range("K1").Value = "7000.0"
$ws.Range("K1").Copy() | Out-Null
$ws.Range("I6").PasteSpecial(-4163) | Out-Null

$ws.Range("K1").Value = "6000.0"
$ws.Range("K1").Copy() | Out-Null
$ws.Range("I7").PasteSpecial(-4163) | Out-Null

$ws.Range("K1").Clear() | Out-Null

# ---------------------------------------------------------------------------
# 5. Misc sheet/view metadata touched by the diff.
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 22
$ws.Range("J7").Select() | Out-Null
